# Adds a "Name" column (ticker full name) as new column C on every
# worksheet, shifting the existing Price Change / Total Dividend /
# Dividend Yield / Unadjusted Total / Adjusted Total columns from C:G
# to D:H.

$wb = $excel.ActiveWorkbook

# Ticker -> Name lookup (shared across all sheets)
$names = @{
    "TLT"    = "iShares 20+ Year Treasury Bond ETF"
    "^GSPC"  = "S&P 500"
    "DIVO"   = "Amplify CWP Enhanced Dividend Income ETF"
    "FTHI"   = "First Trust BuyWrite Income ETF"
    "SPYI"   = "Neos S&P 500(R) High Income ETF"
    "RDVI"   = "FT Vest Rising Dividend Achievers Target Income ETF"
    "JEPI"   = "JPMorgan Equity Premium Income ETF"
    "ISPY"   = "ProShares S&P 500 High Income ETF"

    "EMLP"   = "First Trust North American Energy Infrastructure Fund"
    "MLPR"   = "ETRACS Quarterly Pay 1.5X Leveraged Alerian MLP Index ETN"
    "AMLP"   = "Alerian MLP ETF"

    "ACWI"   = "iShares MSCI ACWI ETF"
    "ACWV"   = "iShares MSCI Global Min Vol Factor ETF"
    "URTH"   = "iShares MSCI World ETF"
    "FLPSX"  = "Fidelity Low-Priced Stock"
    "GVALX"  = "Gotham Large Value Institutional"
    "VMVFX"  = "Vanguard Global Minimum Volatility Inv"
    "CCGIX"  = "Baird Chautauqua Global Growth Instl"
    "SGENX"  = "First Eagle Global A"
    "WGNIX"  = "Pabrai Wagons Institutional"
    "BRK-A"  = "Berkshire Hathaway Inc."
    "MKL"    = "Markel Group Inc."

    "EWJV"   = "iShares MSCI Japan Value ETF"
    "VWO"    = "Vanguard Emerging Markets Stock Index Fund"
    "HEDJ"   = "WisdomTree Europe Hedged Equity Fund"
    "FNDE"   = "Schwab Fundamental Emerging Markets Equity ETF"
    "ADVE"   = "Matthews Asia Dividend Active ETF"
    "EEM"    = "iShares MSCI Emerging Markets ETF"

    "GVLU"   = "Gotham 1000 Value ETF"
    "GSPY"   = "Gotham Enhanced 500 ETF"
    "RSP"    = "Invesco S&P 500 Equal Weight ETF"
    "IWN"    = "iShares Russell 2000 Value ETF"
    "GINDX"  = "Gotham Index Plus Institutional"
    "GSPFX"  = "Gotham Enhanced S&P 500 Index Instl"
    "SEQUX"  = "Sequoia"
    "FEVAX"  = "First Eagle US Value A"
    "FCNTX"  = "Fidelity Contrafund"
    "FMAGX"  = "Fidelity Magellan"
    "AKRIX"  = "Akre Focus Instl"
}

for ($si = 1; $si -le $wb.Worksheets.Count; $si++) {
    $ws = $wb.Worksheets.Item($si)

    # Insert a new blank column before the current column C; this
    # shifts the old C:G (Price Change..Adjusted Total) to D:H and
    # copies the header style/formatting along with it.
    $ws.Columns.Item(3).Insert()

    $ws.Cells.Item(1, 3).Value = "Name"

    $lastRow = $ws.Cells.Item(1, 1).End(-4121).Row # xlDown

    for ($r = 2; $r -le $lastRow; $r++) {
        $ticker = $ws.Cells.Item($r, 2).Value2
        if ($names.ContainsKey($ticker)) {
            $ws.Cells.Item($r, 3).Value = $names[$ticker]
        }
    }
}
